$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.886.82"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.533.89"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.77"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.09"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +8.97%  "
$ws.Range("E6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E7").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.80"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.96%  "
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.73"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("E12").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.930.13"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.546.25"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.12"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.06%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.870"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.955.32"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.15"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0983"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.49"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.58"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.12"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.08"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -7.06%  "
$ws.Range("E26").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E27").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +9.14%  "
$ws.Range("E28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.22"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.93"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.33%  "
$ws.Range("E30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.64"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.11"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0795"
$ws.Range("D34").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E34").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.18%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.26"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.23"
$ws.Range("D37").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("E38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.25"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.00%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.119"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.46"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E41").ClearFormats()

$ws.Range("B42").Value = "ApeXProtocol"

$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.09"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.38%  "
$ws.Range("E42").ClearFormats()

$ws.Range("B43").Value = "RenderToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.87"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0304"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("E44").ClearFormats()

$ws.Range("B45").Value = "FirstDigitalUSD"

$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E45").ClearFormats()

$ws.Range("B46").Value = "Maker"

$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.047.03"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("E46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.08"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.98"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.788.49"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.192"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("E50").ClearFormats()

$ws.Range("B51").Value = "Aave"

$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.69"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.36%  "
$ws.Range("E51").ClearFormats()
